$wb = $excel.ActiveWorkbook

# Use an existing cell that already carries the "bordered, left/top, wrap" style
# as a format donor so the new sheet reuses the same cellXfs entry instead of
# creating a brand new one.
$ws2 = $wb.Worksheets.Item("Test cases")
$styleDonor = $ws2.Range("C2")

# Add the new "Defects" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Defects"

# Apply the same formatting (border + left/top + wrap) used elsewhere in the
# workbook to the whole data range before filling in values.
$styleDonor.Copy()
$ws.Range("A1:D4").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$ws.Range("A1").Value = "SrNo"
$ws.Range("B1").Value = "Defect Description"
$ws.Range("C1").Value = "Defect Steps"
$ws.Range("D1").Value = "Actual condition"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "December month page shows blank page"
$ws.Range("C2").Value = "1. Open rb shoe site`n2. Click December month link`n3. Verify that if no shoe present it should show Coming Soon text"
$ws.Range("D2").Value = "December page shows blank without any message"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "November month page does not show shoe image"
$ws.Range("C3").Value = "1. Open rb shoe site`n2. Click November month link`n3. Verify that page should show image of shoe"
$ws.Range("D3").Value = "November month page does not show shoe image"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Month page: Description tag alignment not proper"
$ws.Range("C4").Value = "1. Open rb shoe site`n2. Open any month page`n3. Verify Description tag alignment"
$ws.Range("D4").Value = "Description tag alignment is not proper"

# Row heights
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 45

# Column widths for B:D (closest achievable value to the target 41.7109375
# "characters" given this host's internal column-width quantization).
$ws.Range("B1:D1").ColumnWidth = 40.8

# Make the new "Defects" sheet the selected/active tab.
$ws.Activate()

Write-Host "Defects sheet created"
